$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row numbers in column A (test-case IDs) ---
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# --- New text, entered in the order that reproduces the original shared-string order ---
$ws.Range("C4").Value = "Jasmine"
$ws.Range("D4").Value = "Valid email & invalid password"
$ws.Range("B7").Value = "123@gmail.com"
$ws.Range("D3").Value = "Null email & null password"
$ws.Range("D5").Value = "Valid email & null password"
$ws.Range("D6").Value = "Null email & password"
$ws.Range("D7").Value = "Invalid email & password"
$ws.Range("D8").Value = "Invalid email & null password"
$ws.Range("B5").Value = "jasmine.liu012005@gmail.com"
$ws.Range("B8").Value = "345@mail.com"
$ws.Range("B10").Value = "123@mail"
$ws.Range("D10").Value = "You've subscribed to a GrabOne newsletter, but you haven't registered. `nPlease register here."

# --- Cells that reuse already-existing text ---
$ws.Range("B4").Value = "lovelydoudou0814@gmail.com"
$ws.Range("C6").Value = "Jasmine"
$ws.Range("C7").Value = "Jasmine0814"

# --- New hyperlinks (email addresses become mailto: links, like the existing B2 one) ---
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:lovelydoudou0814@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:jasmine.liu012005@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:123@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:345@mail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:123@mail") | Out-Null

# --- Row 8 grows a proper bottom border to match the rest of the table ---
$ws.Range("A8:D8").Borders.LineStyle = 1

# --- The note in row 10 wraps and the row grows to fit it ---
$ws.Range("D10").WrapText = $true
$ws.Rows(10).RowHeight = 45

# --- Column D widens to fit the new "Note" text ---
$ws.Columns("D").ColumnWidth = 36.66

# --- Restore the cursor position as saved in the workbook ---
$ws.Range("D15").Select() | Out-Null
